{"js": "// Re-apply the \"Couldn't find / bookmark / REF bookmark1\" bookmark-link\n// fixture changes (m2doc #476, POI 4.1.0 -> 5.2.3 regeneration):\n//  - keep the \"Couldn't find the 'link' variable\" run bold\n//  - keep the \"a reference to bookmark1\" REF-field result run bold\n//  - refresh the \"bookmark1\" bookmark (delete + re-insert on the same\n//    range) so its start/end markers are regenerated, matching the\n//    bookmark id refresh seen in the diff.\nconst doc = context.document;\nconst body = doc.body;\n\n// 1) \"Couldn't find the 'link' variable\" stays bold.\nconst errorResults = body.search(\"Couldn't find the 'link' variable\", { matchCase: true });\nerrorResults.load(\"items\");\nawait context.sync();\nif (errorResults.items.length > 0) {\n  errorResults.items[0].font.bold = true;\n}\n\n// 2) The REF field's cached result text (\"a reference to bookmark1\") stays bold.\nconst refResults = body.search(\"a reference to bookmark1\", { matchCase: true });\nrefResults.load(\"items\");\nawait context.sync();\nif (refResults.items.length > 0) {\n  refResults.items[0].font.bold = true;\n}\n\n// 3) Recreate bookmark \"bookmark1\" in place so it gets a fresh identity,\n//    like the regenerated bookmarkStart/bookmarkEnd ids in the diff.\nconst bookmarkRange = doc.getBookmarkRange(\"bookmark1\");\ndoc.deleteBookmark(\"bookmark1\");\nbookmarkRange.insertBookmark(\"bookmark1\");\n\nawait context.sync();\n", "ps1": "# Re-apply the \"Couldn't find / bookmark / REF bookmark1\" bookmark-link\n# fixture changes (m2doc #476, POI 4.1.0 -> 5.2.3 regeneration):\n#  - keep the \"Couldn't find the 'link' variable\" run bold\n#  - keep the \"a reference to bookmark1\" REF-field result run bold\n#  - refresh the \"bookmark1\" bookmark (delete + re-insert on the same\n#    range) so its start/end markers are regenerated, matching the\n#    bookmark id refresh seen in the diff.\n$d = $word.ActiveDocument\n\n# 1) \"Couldn't find the 'link' variable\" stays bold.\n$errRange = $d.Content\n$found1 = $errRange.Find.Execute(\"Couldn't find the 'link' variable\")\nif ($found1) {\n    $errRun = $d.Range($errRange.Start, $errRange.End)\n    $errRun.Font.Bold = 1\n}\n\n# 2) The REF field's cached result text (\"a reference to bookmark1\") stays bold.\nforeach ($fld in $d.Fields) {\n    if ($fld.Code.Text -match \"REF bookmark1\") {\n        $fld.Result.Font.Bold = 1\n    }\n}\n\n# 3) Recreate bookmark \"bookmark1\" in place so it gets a fresh identity,\n#    like the regenerated bookmarkStart/bookmarkEnd ids in the diff.\n$bm = $d.Bookmarks(\"bookmark1\")\n$bmRange = $bm.Range\n$bm.Delete()\n$d.Bookmarks.Add(\"bookmark1\", $bmRange)\n"}
